# GCF_File_Usage.xlsx - append new usage-log rows (65711.. -> rows 64-77)
# to the "Data" sheet, extending the used range from A1:O63 to A1:O77.

$wb = $excel.ActiveWorkbook

# Locate worksheet "Data" (falls back to the active sheet if not found).
$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Data") {
        $ws = $sheet
        break
    }
}
if ($null -eq $ws) {
    $ws = $wb.ActiveSheet
}

# New rows to append (columns A:O), in order, matching the commit's data.
$newRows = @(
    @(45711.627835648149, 8, 6, 217, 392, 368, 388, 2681, 388, 1216, 119, 304, 30, 3129, 4301),
    @(45712.930486111109, 8, 6, 223, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3157, 4344),
    @(45712.93478009259,  8, 6, 223, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3157, 4344),
    @(45712.935312499998, 8, 6, 223, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3157, 4344),
    @(45713.945706018516, 8, 6, 223, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3157, 4408),
    @(45713.94672453704,  8, 6, 223, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3157, 4408),
    @(45713.948530092595, 8, 6, 223, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3157, 4408),
    @(45714.709872685184, 8, 6, 244, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3232, 4408),
    @(45714.712337962963, 8, 6, 244, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3232, 4408),
    @(45714.717083333337, 8, 6, 244, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3232, 4408),
    @(45714.725185185183, 8, 6, 244, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3232, 4408),
    @(45714.731226851851, 8, 6, 244, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3232, 4408),
    @(45714.932314814818, 8, 6, 229, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3181, 4451),
    @(45714.937083333331, 8, 6, 229, 393, 369, 388, 2681, 388, 1216, 119, 304, 30, 3181, 4451)
)

# The sheet currently uses rows 1 (header) through 63 (last data row), so the
# new block starts at row 64 and extends the dimension to A1:O77.
$startRow = 64
$numRows  = $newRows.Count
$numCols  = 15

$values = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $values[$i, $j] = $newRows[$i][$j]
    }
}

$endRow = $startRow + $numRows - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $numCols))
$targetRange.Value2 = $values
